$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column AR ("valor") holds 25 for every data row (rows 2-51); update it to 16
# (re-labelling this "Q" letter sample set, per commit "letras pablo con etiqueta").
$ws.Range("AR2:AR51").Value = 16

# Reflect the cell that was selected/active when the workbook was last saved.
$ws.Range("AS49").Select()
